# Update "F" column (想去人数 / interest counts) values across sheets to
# match the newly generated data output ("Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2145
$ws1.Range("F10").Value = 43
$ws1.Range("F11").Value = 2549
$ws1.Range("F12").Value = 1615
$ws1.Range("F13").Value = 1593
$ws1.Range("F15").Value = 262
$ws1.Range("F16").Value = 640
$ws1.Range("F17").Value = 815
$ws1.Range("F18").Value = 99
$ws1.Range("F19").Value = 326
$ws1.Range("F24").Value = 5468
$ws1.Range("F26").Value = 807
$ws1.Range("F27").Value = 100
$ws1.Range("F31").Value = 229
$ws1.Range("F32").Value = 43
$ws1.Range("F33").Value = 1056
$ws1.Range("F34").Value = 785
$ws1.Range("F38").Value = 412
$ws1.Range("F44").Value = 87

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 790

# --- Sheet: 全部类型 (All types, aggregate of the other sheets) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2145
$ws4.Range("F14").Value = 43
$ws4.Range("F15").Value = 2549
$ws4.Range("F16").Value = 1615
$ws4.Range("F17").Value = 1593
$ws4.Range("F19").Value = 262
$ws4.Range("F20").Value = 640
$ws4.Range("F22").Value = 815
$ws4.Range("F23").Value = 99
$ws4.Range("F24").Value = 326
$ws4.Range("F28").Value = 5468
$ws4.Range("F30").Value = 807
$ws4.Range("F31").Value = 100
$ws4.Range("F35").Value = 229
$ws4.Range("F36").Value = 43
$ws4.Range("F37").Value = 1056
$ws4.Range("F38").Value = 785
$ws4.Range("F40").Value = 412
$ws4.Range("F46").Value = 87
